$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 671
$ws.Range("I33").Value = 634
$ws.Range("J33").Value = 819
$ws.Range("K33").Value = 634
$ws.Range("L33").Value = 819
$ws.Range("M33").Value = -405
$ws.Range("N33").Value = -1277
$ws.Range("H98").Value = 213104.45
$ws.Range("I98").Value = 1340.55
$ws.Range("K98").Value = 1340.55
$ws.Range("M98").Value = 157.45
$ws.Range("H113").Value = 6749.4614
$ws.Range("J113").Value = 7370.5
$ws.Range("L113").Value = 7370.5
$ws.Range("N113").Value = -13878.5
$ws.Range("H122").Value = 213104.45
$ws.Range("I122").Value = 1340.55
$ws.Range("K122").Value = 4021.65
$ws.Range("M122").Value = -1571.65
$ws.Range("H132").Value = 1094.8096
$ws.Range("I132").Value = 949.55
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 2848.65
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -318.6499999999996
$ws.Range("N132").Value = -17060
$ws.Range("H138").Value = 2957.0234
$ws.Range("I138").Value = 1715.3182
$ws.Range("J138").Value = 3390.635
$ws.Range("K138").Value = 5145.9546
$ws.Range("L138").Value = 10171.905
$ws.Range("M138").Value = -5.954600000000028
$ws.Range("N138").Value = -20451.905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10755357
$ws.Range("I74").Value = 15153521
$ws.Range("K74").Value = 15153521
$ws.Range("M74").Value = -15152647
$ws.Range("H77").Value = 10755357
$ws.Range("I77").Value = 15153521
$ws.Range("K77").Value = 75767605
$ws.Range("M77").Value = -75763237
$ws.Range("H97").Value = 1763.625
$ws.Range("I97").Value = 1844.2858
$ws.Range("K97").Value = 1844.2858
$ws.Range("M97").Value = -1348.2858
$ws.Range("H122").Value = 5099.8237
$ws.Range("I122").Value = 4909.5557
$ws.Range("J122").Value = 5313.875
$ws.Range("K122").Value = 14728.6671
$ws.Range("L122").Value = 15941.625
$ws.Range("M122").Value = -12278.6671
$ws.Range("N122").Value = -20841.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1346
$ws.Range("H29").Value = 9963
$ws.Range("I29").Value = 9963
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 9963
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -9674
$ws.Range("N29").ClearContents()
$ws.Range("H64").Value = 2000
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 2000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -1775
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 2000
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 2000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -1220
$ws.Range("N67").ClearContents()
$ws.Range("H94").Value = 4282.75
$ws.Range("I94").Value = 4125.3335
$ws.Range("K94").Value = 4125.3335
$ws.Range("M94").Value = -3674.3335
$ws.Range("H99").Value = 1700
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H128").Value = 17969
$ws.Range("I128").Value = 17969
$ws.Range("K128").Value = 53907
$ws.Range("M128").Value = -51417

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 302.66666
$ws.Range("I7").Value = 59.5
$ws.Range("K7").Value = 59.5
$ws.Range("M7").Value = 53.5
$ws.Range("H31").Value = 31998.676
$ws.Range("I31").Value = 1959.5416
$ws.Range("K31").Value = 1959.5416
$ws.Range("M31").Value = -1664.5416
$ws.Range("H34").Value = 31998.676
$ws.Range("I34").Value = 1959.5416
$ws.Range("K34").Value = 1959.5416
$ws.Range("M34").Value = -1757.5416
$ws.Range("H99").Value = 2107.8
$ws.Range("I99").Value = 1950
$ws.Range("K99").Value = 1950
$ws.Range("M99").Value = -452
$ws.Range("H126").Value = 2107.8
$ws.Range("I126").Value = 1950
$ws.Range("K126").Value = 5850
$ws.Range("M126").Value = -3380
$ws.Range("H135").Value = 69418.39999999999
$ws.Range("J135").Value = 69418.39999999999
$ws.Range("L135").Value = 69418.39999999999
$ws.Range("N135").Value = -79558.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 120789.8
$ws.Range("I128").Value = 120789.8
$ws.Range("K128").Value = 362369.4
$ws.Range("M128").Value = -357389.4
$ws.Range("H131").Value = 7100362.5
$ws.Range("J131").Value = 4744200
$ws.Range("L131").Value = 14232600
$ws.Range("N131").Value = -14242680

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5416
$ws.Range("I7").Value = 4077.6
$ws.Range("J7").Value = 9598.5
$ws.Range("K7").Value = 4077.6
$ws.Range("L7").Value = 9598.5
$ws.Range("M7").Value = -3965.6
$ws.Range("N7").Value = -9822.5
$ws.Range("H40").Value = 10758.608
$ws.Range("I40").Value = 11117.571
$ws.Range("J40").Value = 10200.223
$ws.Range("K40").Value = 11117.571
$ws.Range("L40").Value = 10200.223
$ws.Range("M40").Value = -10981.571
$ws.Range("N40").Value = -10472.223
$ws.Range("H68").Value = 2034.3636
$ws.Range("I68").Value = 877.2
$ws.Range("J68").Value = 2998.6667
$ws.Range("K68").Value = 877.2
$ws.Range("L68").Value = 2998.6667
$ws.Range("M68").Value = -128.2
$ws.Range("N68").Value = -4496.6667
$ws.Range("H71").Value = 2034.3636
$ws.Range("I71").Value = 877.2
$ws.Range("J71").Value = 2998.6667
$ws.Range("K71").Value = 4386
$ws.Range("L71").Value = 14993.3335
$ws.Range("M71").Value = -642
$ws.Range("N71").Value = -22481.3335
$ws.Range("H93").Value = 1711.75
$ws.Range("I93").Value = 1644
$ws.Range("K93").Value = 1644
$ws.Range("M93").Value = -396
$ws.Range("H94").Value = 10000
$ws.Range("J94").Value = 10000
$ws.Range("L94").Value = 10000
$ws.Range("N94").Value = -11352
$ws.Range("H126").Value = 5416
$ws.Range("I126").Value = 4077.6
$ws.Range("J126").Value = 9598.5
$ws.Range("K126").Value = 12232.8
$ws.Range("L126").Value = 28795.5
$ws.Range("M126").Value = -9762.799999999999
$ws.Range("N126").Value = -33735.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3622.4443
$ws.Range("I81").Value = 3950
$ws.Range("J81").Value = 3360.4
$ws.Range("K81").Value = 7900
$ws.Range("L81").Value = 6720.8
$ws.Range("M81").Value = -6839
$ws.Range("N81").Value = -8842.799999999999
$ws.Range("H84").Value = 3622.4443
$ws.Range("I84").Value = 3950
$ws.Range("J84").Value = 3360.4
$ws.Range("K84").Value = 39500
$ws.Range("L84").Value = 33604
$ws.Range("M84").Value = -34196
$ws.Range("N84").Value = -44212
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H122").Value = 1902.7368
$ws.Range("I122").Value = 1256.875
$ws.Range("J122").Value = 5347.3335
$ws.Range("K122").Value = 3770.625
$ws.Range("L122").Value = 16042.0005
$ws.Range("M122").Value = -1320.625
$ws.Range("N122").Value = -20942.0005
$ws.Range("H132").Value = 3722.2144
$ws.Range("I132").Value = 2092.0833
$ws.Range("K132").Value = 6276.249899999999
$ws.Range("M132").Value = -3746.249899999999
